# Form the consolidated report: populate the "Absent" (column H) figures
# for each attendance date row so that a student is marked absent (1)
# whenever they were not marked present ("Real" / column E = 0), and
# present (0) otherwise. This both fixes rows that previously held 0
# by mistake and fills in rows that were left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value  = 1
$ws.Range("H6").Value  = 0
$ws.Range("H10").Value = 1
$ws.Range("H13").Value = 0
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 0
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 0
